$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for the two date groups (2-5 and 6-9) need to swap their
# Fecha / Volumen / Precio mínimo / Precio máximo / Precio promedio
# ponderado / Origen / Precio $/Kg values, while everything else
# (Mercado, Región, Codreg, Tipo, Producto, Categoría, Variedad, Calidad,
# Unidad de comercialización, Kg / unidad) stays put.

$rowsTop = 2..5
$rowsBottom = 6..9
$cols = @("D", "M", "N", "O", "P", "R", "S")

for ($i = 0; $i -lt $rowsTop.Count; $i++) {
    $rTop = $rowsTop[$i]
    $rBottom = $rowsBottom[$i]

    foreach ($col in $cols) {
        $topCell = $ws.Range("$col$rTop")
        $bottomCell = $ws.Range("$col$rBottom")

        $topValue = $topCell.Value2
        $bottomValue = $bottomCell.Value2

        $topCell.Value2 = $bottomValue
        $bottomCell.Value2 = $topValue
    }
}

$wb.Save()
